$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "removed fake low input data" -----------------------------------------
# The sheet held a dilution series (lasAHL_nM vs MEFL_geomean). The very low,
# fabricated concentration points (0.01/4.1, 1E-3/4, 1E-4/3.9, 1E-5/3.8) are
# replaced by a single real "zero input" control point, and the trailing
# rows are dropped so the used range shrinks from A1:B17 to A1:B14.

# Row 13 becomes the zero-concentration control point.
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 3.76

# Row 14 is kept as an empty (but still styled) row...
$ws.Range("A14").Value = $null
$ws.Range("B14").Value = $null

# ...and old rows 15-17 (1E-4, 1E-5, and the trailing blank row) are removed
# outright, shifting everything below them up.
$ws.Range("A15:B17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null

# --- restore the view/selection recorded in the saved file -----------------
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G10").Select() | Out-Null
